$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (manual_test_join_1)
$ws.Range("E2").Value = 2

# Row 3 (high_level_test_join_1)
$ws.Range("D3").Value = 170
$ws.Range("E3").Value = 2

# Row 4 (manual_test_join_2)
$ws.Range("E4").Value = 2

# Row 5 (high_level_test_join_2)
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 145
$ws.Range("E5").Value = 2

# Row 6 (manual_test_join_3)
$ws.Range("E6").Value = 3

# Row 7 (high_level_test_join_3)
$ws.Range("C7").Value = 4
$ws.Range("E7").Value = 2

# Row 8 (manual_test_join_8)
$ws.Range("E8").Value = 5

# Row 9 (high_level_test_join_8)
$ws.Range("E9").Value = 4

# Row 10 (manual_test_agg_1)
$ws.Range("E10").Value = 8

# Row 12 (manual_test_agg_6)
$ws.Range("E12").Value = 9

# Row 13 (high_level_test_agg_6)
$ws.Range("E13").Value = 3

# Row 14 (manual_test_distinct_2)
$ws.Range("E14").Value = 5

# Row 16 (manual_test_distinct_4)
$ws.Range("E16").Value = 6

# Row 17 (high_level_test_distinct_4)
$ws.Range("E17").Value = 2

# Row 18 (manual_test_where_1)
$ws.Range("E18").Value = 7

# Row 19 (high_level_test_where_1)
$ws.Range("E19").Value = 3

# Row 20 (manual_test_where_having_1)
$ws.Range("E20").Value = 10

# Row 21 (high_level_test_where_having_1)
$ws.Range("E21").Value = 2

# Row 22 (manual_test_large_query_1)
$ws.Range("E22").Value = 8

# Row 23 (high_level_test_large_query_1)
$ws.Range("E23").Value = 5

# Row 24 (manual_test_large_query_3)
$ws.Range("E24").Value = 13

# Row 25 (high_level_test_large_query_3)
$ws.Range("E25").Value = 5

# Row 26 (manual_test_one_cell_3)
$ws.Range("E26").Value = 10

# Row 27 (high_level_test_one_cell_3)
$ws.Range("E27").Value = 2

# Row 28 (manual_test_one_cell_5)
$ws.Range("E28").Value = 11

# Row 29 (high_level_test_one_cell_5)
$ws.Range("E29").Value = 3

# Row 30 (manual_test_mixed_data_1)
$ws.Range("E30").Value = 6

# Row 31 (high_level_test_mixed_data_1)
$ws.Range("E31").Value = 2

# Row 32 (manual_test_mixed_data_3)
$ws.Range("E32").Value = 5

# Row 33 (high_level_test_mixed_data_3)
$ws.Range("E33").Value = 3
